$wb = $excel.ActiveWorkbook

# --- Sheet 1: Significant Components ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(1,2).Value = "Factor"
$ws.Cells.Item(1,3).Value = "Sig Components"
$ws.Cells.Item(2,1).Value = "F2: 1"
$ws.Cells.Item(2,2).Value = "F2: 1"
$ws.Cells.Item(2,3).Value = "['QNOHLTH' 'QHISPC' 'PPUNIT' 'QSERV' 'QESL' 'QEDLESHI' 'QEXTRCT' 'QFHH'`n 'PERCAP']"
$ws.Cells.Item(3,1).Value = "F2: 2"
$ws.Cells.Item(3,2).Value = "F2: 2"
$ws.Cells.Item(3,3).Value = "['QRICH' 'MDHSEVAL' 'PERCAP']"
$ws.Cells.Item(4,1).Value = "F2: 3"
$ws.Cells.Item(4,2).Value = "F2: 3"
$ws.Cells.Item(4,3).Value = "['QRENTER' 'QNOAUTO' 'QPOVTY']"
$ws.Cells.Item(5,1).Value = "F2: 4"
$ws.Cells.Item(5,2).Value = "F2: 4"
$ws.Cells.Item(5,3).Value = "['MEDAGE' 'QSSBEN' 'QAGEDEP']"
$ws.Cells.Item(6,1).Value = "F2: 5"
$ws.Cells.Item(6,2).Value = "F2: 5"
$ws.Cells.Item(6,3).Value = "['QAGEDEP' 'QFEMLBR' 'QFEMALE']"
$ws.Cells.Item(1,2).Copy() | Out-Null
$ws.Cells.Item(6,1).PasteSpecial(-4122) | Out-Null

# --- Sheet 2: Loading Factors ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(1,2).Value = "F2: 1"
$ws.Cells.Item(1,3).Value = "F2: 2"
$ws.Cells.Item(1,4).Value = "F2: 3"
$ws.Cells.Item(1,5).Value = "F2: 4"
$ws.Cells.Item(1,6).Value = "F2: 5"
$ws.Cells.Item(2,1).Value = "QNOHLTH"
$ws.Cells.Item(2,2).Value = 0.6641941122002464
$ws.Cells.Item(2,3).Value = 0.4111007544969877
$ws.Cells.Item(2,4).Value = 0.327159092345672
$ws.Cells.Item(2,5).Value = -0.08880591641510853
$ws.Cells.Item(2,6).Value = -0.1206077732995754
$ws.Cells.Item(3,1).Value = "QHISPC"
$ws.Cells.Item(3,2).Value = 0.817280070824076
$ws.Cells.Item(3,3).Value = 0.3568107832167007
$ws.Cells.Item(3,4).Value = 0.1430522455328176
$ws.Cells.Item(3,5).Value = -0.1035311703573437
$ws.Cells.Item(3,6).Value = -0.1284245812002146
$ws.Cells.Item(4,1).Value = "PPUNIT"
$ws.Cells.Item(4,2).Value = 0.7213533591425673
$ws.Cells.Item(4,3).Value = 0.03768695314523186
$ws.Cells.Item(4,4).Value = -0.340714441132671
$ws.Cells.Item(4,5).Value = -0.0929331490039619
$ws.Cells.Item(4,6).Value = 0.07082527372821094
$ws.Cells.Item(5,1).Value = "QSERV"
$ws.Cells.Item(5,2).Value = 0.517056355052926
$ws.Cells.Item(5,3).Value = 0.3832057679581961
$ws.Cells.Item(5,4).Value = 0.311558311947261
$ws.Cells.Item(5,5).Value = -0.1307287334302756
$ws.Cells.Item(5,6).Value = -0.1013539558858636
$ws.Cells.Item(6,1).Value = "QESL"
$ws.Cells.Item(6,2).Value = 0.7916387764699077
$ws.Cells.Item(6,3).Value = 0.1578130985349176
$ws.Cells.Item(6,4).Value = 0.2139992488008075
$ws.Cells.Item(6,5).Value = -0.0314487937301029
$ws.Cells.Item(6,6).Value = -0.2060743019298396
$ws.Cells.Item(7,1).Value = "QEDLESHI"
$ws.Cells.Item(7,2).Value = 0.8729342608933685
$ws.Cells.Item(7,3).Value = 0.2064110063017868
$ws.Cells.Item(7,4).Value = 0.2004013888253722
$ws.Cells.Item(7,5).Value = -0.01618995688075913
$ws.Cells.Item(7,6).Value = -0.1067088317849588
$ws.Cells.Item(8,1).Value = "QEXTRCT"
$ws.Cells.Item(8,2).Value = 0.7830485250109201
$ws.Cells.Item(8,3).Value = 0.1226948168972316
$ws.Cells.Item(8,4).Value = 0.09276530213712816
$ws.Cells.Item(8,5).Value = -0.01971134291552629
$ws.Cells.Item(8,6).Value = -0.2082444669752325
$ws.Cells.Item(9,1).Value = "QFHH"
$ws.Cells.Item(9,2).Value = 0.5607864131566372
$ws.Cells.Item(9,3).Value = 0.2738331457858463
$ws.Cells.Item(9,4).Value = 0.05550829987809005
$ws.Cells.Item(9,5).Value = -0.06567204376371692
$ws.Cells.Item(9,6).Value = 0.2456792748662655
$ws.Cells.Item(10,1).Value = "QRICH"
$ws.Cells.Item(10,2).Value = 0.2058450521114599
$ws.Cells.Item(10,3).Value = 0.8483747168011316
$ws.Cells.Item(10,4).Value = 0.3541943659819356
$ws.Cells.Item(10,5).Value = -0.1749190839245255
$ws.Cells.Item(10,6).Value = -0.02690806984672741
$ws.Cells.Item(11,1).Value = "MDHSEVAL"
$ws.Cells.Item(11,2).Value = 0.3746877273915546
$ws.Cells.Item(11,3).Value = 0.680310731367022
$ws.Cells.Item(11,4).Value = -0.1250235556616157
$ws.Cells.Item(11,5).Value = 0.00295119629720828
$ws.Cells.Item(11,6).Value = -0.01357761256245511
$ws.Cells.Item(12,1).Value = "PERCAP"
$ws.Cells.Item(12,2).Value = 0.472078660569859
$ws.Cells.Item(12,3).Value = 0.7110281175952148
$ws.Cells.Item(12,4).Value = 0.2599191305489495
$ws.Cells.Item(12,5).Value = -0.2188984120719565
$ws.Cells.Item(12,6).Value = 0.04385980555370568
$ws.Cells.Item(13,1).Value = "QRENTER"
$ws.Cells.Item(13,2).Value = -0.0169865811594124
$ws.Cells.Item(13,3).Value = 0.1782040321943912
$ws.Cells.Item(13,4).Value = 0.7623798549014968
$ws.Cells.Item(13,5).Value = -0.4175899392239148
$ws.Cells.Item(13,6).Value = -0.1129122996727554
$ws.Cells.Item(14,1).Value = "QNOAUTO"
$ws.Cells.Item(14,2).Value = 0.1056780863584004
$ws.Cells.Item(14,3).Value = 0.05350761822061592
$ws.Cells.Item(14,4).Value = 0.6975393916527455
$ws.Cells.Item(14,5).Value = -0.04914871524614545
$ws.Cells.Item(14,6).Value = 0.003119068169388665
$ws.Cells.Item(15,1).Value = "QPOVTY"
$ws.Cells.Item(15,2).Value = 0.2992888396138256
$ws.Cells.Item(15,3).Value = 0.148607011408278
$ws.Cells.Item(15,4).Value = 0.5737043965764143
$ws.Cells.Item(15,5).Value = -0.2807685012847285
$ws.Cells.Item(15,6).Value = 0.09516280552966946
$ws.Cells.Item(16,1).Value = "MEDAGE"
$ws.Cells.Item(16,2).Value = -0.3193244214567408
$ws.Cells.Item(16,3).Value = -0.2200222480200467
$ws.Cells.Item(16,4).Value = -0.3643133272133273
$ws.Cells.Item(16,5).Value = 0.7585327102945554
$ws.Cells.Item(16,6).Value = -0.04004000901821723
$ws.Cells.Item(17,1).Value = "QSSBEN"
$ws.Cells.Item(17,2).Value = 0.005424293531374609
$ws.Cells.Item(17,3).Value = -0.02898540048679379
$ws.Cells.Item(17,4).Value = -0.1541266274192051
$ws.Cells.Item(17,5).Value = 0.829857786771077
$ws.Cells.Item(17,6).Value = 0.08691606186838186
$ws.Cells.Item(18,1).Value = "QAGEDEP"
$ws.Cells.Item(18,2).Value = -0.01574198855808457
$ws.Cells.Item(18,3).Value = -0.1414286087818206
$ws.Cells.Item(18,4).Value = -0.1080367958173994
$ws.Cells.Item(18,5).Value = 0.6951599658178635
$ws.Cells.Item(18,6).Value = 0.5982193179122143
$ws.Cells.Item(19,1).Value = "QFEMLBR"
$ws.Cells.Item(19,2).Value = -0.2609508547515848
$ws.Cells.Item(19,3).Value = 0.07812259117665589
$ws.Cells.Item(19,4).Value = 0.00002961441042594403
$ws.Cells.Item(19,5).Value = -0.04558978074484524
$ws.Cells.Item(19,6).Value = 0.7376544298789692
$ws.Cells.Item(20,1).Value = "QFEMALE"
$ws.Cells.Item(20,2).Value = -0.01186351313103661
$ws.Cells.Item(20,3).Value = -0.0748927632368649
$ws.Cells.Item(20,4).Value = 0.008529709894781311
$ws.Cells.Item(20,5).Value = 0.1928293016857023
$ws.Cells.Item(20,6).Value = 0.868720233107825
$ws.Cells.Item(1,2).Copy() | Out-Null
$ws.Cells.Item(1,6).PasteSpecial(-4122) | Out-Null

# --- Sheet 3: All Refactor Variances ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(1,2).Value = "F0: 1"
$ws.Cells.Item(1,3).Value = "F0: 2"
$ws.Cells.Item(1,4).Value = "F0: 3"
$ws.Cells.Item(1,5).Value = "F0: 4"
$ws.Cells.Item(1,6).Value = "F0: 5"
$ws.Cells.Item(1,7).Value = "F0: 6"
$ws.Cells.Item(1,8).Value = "F0: 7"
$ws.Cells.Item(1,9).Value = "F1: 1"
$ws.Cells.Item(1,10).Value = "F1: 2"
$ws.Cells.Item(1,11).Value = "F1: 3"
$ws.Cells.Item(1,12).Value = "F1: 4"
$ws.Cells.Item(1,13).Value = "F1: 5"
$ws.Cells.Item(1,14).Value = "F2: 1"
$ws.Cells.Item(1,15).Value = "F2: 2"
$ws.Cells.Item(1,16).Value = "F2: 3"
$ws.Cells.Item(1,17).Value = "F2: 4"
$ws.Cells.Item(1,18).Value = "F2: 5"
$ws.Cells.Item(2,1).Value = "SS Loadings"
$ws.Cells.Item(2,2).Value = 4.931282064630651
$ws.Cells.Item(2,3).Value = 2.797028014906584
$ws.Cells.Item(2,4).Value = 2.452495185487878
$ws.Cells.Item(2,5).Value = 2.338350485850176
$ws.Cells.Item(2,6).Value = 1.91555137023628
$ws.Cells.Item(2,7).Value = 1.357471911060389
$ws.Cells.Item(2,8).Value = 1.235122303123883
$ws.Cells.Item(2,9).Value = 4.626008835319829
$ws.Cells.Item(2,10).Value = 3.164090131524766
$ws.Cells.Item(2,11).Value = 2.337676532711858
$ws.Cells.Item(2,12).Value = 2.168885487825472
$ws.Cells.Item(2,13).Value = 1.890316422487597
$ws.Cells.Item(2,14).Value = 4.890284008387104
$ws.Cells.Item(2,15).Value = 2.42792500097139
$ws.Cells.Item(2,16).Value = 2.212082293812811
$ws.Cells.Item(2,17).Value = 2.170978841433984
$ws.Cells.Item(2,18).Value = 1.894385220306634
$ws.Cells.Item(3,1).Value = "Proportion Variance"
$ws.Cells.Item(3,2).Value = 0.1826400764678019
$ws.Cells.Item(3,3).Value = 0.1035936301817253
$ws.Cells.Item(3,4).Value = 0.09083315501806954
$ws.Cells.Item(3,5).Value = 0.08660557355000652
$ws.Cells.Item(3,6).Value = 0.07094634704578816
$ws.Cells.Item(3,7).Value = 0.05027673744668107
$ws.Cells.Item(3,8).Value = 0.04574527048606975
$ws.Cells.Item(3,9).Value = 0.22028613501523
$ws.Cells.Item(3,10).Value = 0.1506709586440365
$ws.Cells.Item(3,11).Value = 0.1113179301291361
$ws.Cells.Item(3,12).Value = 0.1032802613250225
$ws.Cells.Item(3,13).Value = 0.09001506773750462
$ws.Cells.Item(3,14).Value = 0.2573833688624791
$ws.Cells.Item(3,15).Value = 0.1277855263669153
$ws.Cells.Item(3,16).Value = 0.1164253838848848
$ws.Cells.Item(3,17).Value = 0.1142620442859991
$ws.Cells.Item(3,18).Value = 0.0997044852792965
$ws.Cells.Item(4,1).Value = "Cumulative Variance"
$ws.Cells.Item(4,2).Value = 0.1826400764678019
$ws.Cells.Item(4,3).Value = 0.2862337066495272
$ws.Cells.Item(4,4).Value = 0.3770668616675967
$ws.Cells.Item(4,5).Value = 0.4636724352176033
$ws.Cells.Item(4,6).Value = 0.5346187822633914
$ws.Cells.Item(4,7).Value = 0.5848955197100725
$ws.Cells.Item(4,8).Value = 0.6306407901961423
$ws.Cells.Item(4,9).Value = 0.22028613501523
$ws.Cells.Item(4,10).Value = 0.3709570936592664
$ws.Cells.Item(4,11).Value = 0.4822750237884025
$ws.Cells.Item(4,12).Value = 0.585555285113425
$ws.Cells.Item(4,13).Value = 0.6755703528509296
$ws.Cells.Item(4,14).Value = 0.2573833688624791
$ws.Cells.Item(4,15).Value = 0.3851688952293944
$ws.Cells.Item(4,16).Value = 0.5015942791142792
$ws.Cells.Item(4,17).Value = 0.6158563234002784
$ws.Cells.Item(4,18).Value = 0.715560808679575
$ws.Cells.Item(5,1).Value = "Ratio Variance"
$ws.Cells.Item(5,2).Value = 0.2896103127281016
$ws.Cells.Item(5,3).Value = 0.1642672529150954
$ws.Cells.Item(5,4).Value = 0.144033111131011
$ws.Cells.Item(5,5).Value = 0.1373294827996622
$ws.Cells.Item(5,6).Value = 0.112498823654782
$ws.Cells.Item(5,7).Value = 0.07972325645324015
$ws.Cells.Item(5,8).Value = 0.0725377603181076
$ws.Cells.Item(5,9).Value = 0.3260743075619216
$ws.Cells.Item(5,10).Value = 0.223027783869141
$ws.Cells.Item(5,11).Value = 0.1647762215428352
$ws.Cells.Item(5,12).Value = 0.152878617140578
$ws.Cells.Item(5,13).Value = 0.1332430698855241
$ws.Cells.Item(5,14).Value = 0.3596946139873548
$ws.Cells.Item(5,15).Value = 0.1785809463247687
$ws.Cells.Item(5,16).Value = 0.162705087356202
$ws.Cells.Item(5,17).Value = 0.1596818088693915
$ws.Cells.Item(5,18).Value = 0.1393375434622828
$ws.Cells.Item(1,2).Copy() | Out-Null
$ws.Cells.Item(1,13).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,14).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,15).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,16).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,17).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,18).PasteSpecial(-4122) | Out-Null

# --- Sheet 4: Final Variances ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(1,2).Value = "F2: 1"
$ws.Cells.Item(1,3).Value = "F2: 2"
$ws.Cells.Item(1,4).Value = "F2: 3"
$ws.Cells.Item(1,5).Value = "F2: 4"
$ws.Cells.Item(1,6).Value = "F2: 5"
$ws.Cells.Item(2,1).Value = "SS Loadings"
$ws.Cells.Item(2,2).Value = 4.890284008387104
$ws.Cells.Item(2,3).Value = 2.42792500097139
$ws.Cells.Item(2,4).Value = 2.212082293812811
$ws.Cells.Item(2,5).Value = 2.170978841433984
$ws.Cells.Item(2,6).Value = 1.894385220306634
$ws.Cells.Item(3,1).Value = "Proportion Variance"
$ws.Cells.Item(3,2).Value = 0.2573833688624791
$ws.Cells.Item(3,3).Value = 0.1277855263669153
$ws.Cells.Item(3,4).Value = 0.1164253838848848
$ws.Cells.Item(3,5).Value = 0.1142620442859991
$ws.Cells.Item(3,6).Value = 0.0997044852792965
$ws.Cells.Item(4,1).Value = "Cumulative Variance"
$ws.Cells.Item(4,2).Value = 0.2573833688624791
$ws.Cells.Item(4,3).Value = 0.3851688952293944
$ws.Cells.Item(4,4).Value = 0.5015942791142792
$ws.Cells.Item(4,5).Value = 0.6158563234002784
$ws.Cells.Item(4,6).Value = 0.715560808679575
$ws.Cells.Item(5,1).Value = "Ratio Variance"
$ws.Cells.Item(5,2).Value = 0.3596946139873548
$ws.Cells.Item(5,3).Value = 0.1785809463247687
$ws.Cells.Item(5,4).Value = 0.162705087356202
$ws.Cells.Item(5,5).Value = 0.1596818088693915
$ws.Cells.Item(5,6).Value = 0.1393375434622828
$ws.Cells.Item(1,2).Copy() | Out-Null
$ws.Cells.Item(1,6).PasteSpecial(-4122) | Out-Null

# --- Sheet 5: Included and Excluded ---
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(1,2).Value = "include"
$ws.Cells.Item(1,3).Value = "exclude"
$ws.Cells.Item(2,1).Value = 0.0
$ws.Cells.Item(2,2).Value = "[['QNOHLTH', 'QHISPC', 'PPUNIT', 'QSERV', 'QESL', 'QEDLESHI', 'QEXTRCT', 'QFHH', 'PERCAP', 'QRICH', 'MDHSEVAL', 'QRENTER', 'QNOAUTO', 'QPOVTY', 'MEDAGE', 'QSSBEN', 'QAGEDEP', 'QFEMLBR', 'QFEMALE']]"
$ws.Cells.Item(2,3).Value = "[['MDGRENT', 'QASIAN', 'QBLACK', 'QCVLUN', 'QFAM', 'QMOHO', 'QNATIVE', 'QUNOCCHU']]"
